# Auto-generated Excel COM-interop script to apply profit/price updates
# sourced from a scheduled market-data refresh across all Disposal/Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1374.8182
$ws.Range("I15").Value = 1374.8182
$ws.Range("K15").Value = 4124.4546
$ws.Range("M15").Value = -3955.4546

$ws.Range("H62").Value = 85750
$ws.Range("I62").Value = 103166.664
$ws.Range("J62").Value = 78285.71000000001
$ws.Range("K62").Value = 103166.664
$ws.Range("L62").Value = 78285.71000000001
$ws.Range("M62").Value = -102542.664
$ws.Range("N62").Value = -79533.71000000001

$ws.Range("H65").Value = 85750
$ws.Range("I65").Value = 103166.664
$ws.Range("J65").Value = 78285.71000000001
$ws.Range("K65").Value = 515833.32
$ws.Range("L65").Value = 391428.55
$ws.Range("M65").Value = -512713.32
$ws.Range("N65").Value = -397668.55

$ws.Range("H111").Value = 33346.1
$ws.Range("I111").Value = 18959.6
$ws.Range("K111").Value = 56878.8
$ws.Range("M111").Value = -53811.8

$ws.Range("H113").Value = 90913130
$ws.Range("J113").Value = 3999
$ws.Range("L113").Value = 3999
$ws.Range("N113").Value = -10507

$ws.Range("H116").Value = 27837080
$ws.Range("I116").Value = 19321258
$ws.Range("K116").Value = 19321258
$ws.Range("M116").Value = -19317816

$ws.Range("H132").Value = 3951.644
$ws.Range("I132").Value = 4610.2583
$ws.Range("J132").Value = 3222.4644
$ws.Range("K132").Value = 13830.7749
$ws.Range("L132").Value = 9667.393199999999
$ws.Range("M132").Value = -11300.7749
$ws.Range("N132").Value = -14727.3932

$ws.Range("H141").Value = 3516.1765
$ws.Range("I141").Value = 3423.4375
$ws.Range("K141").Value = 10270.3125
$ws.Range("M141").Value = -5090.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 219910.73
$ws.Range("I32").Value = 240759.4
$ws.Range("J32").Value = 999.75
$ws.Range("K32").Value = 240759.4
$ws.Range("L32").Value = 999.75
$ws.Range("M32").Value = -240472.4
$ws.Range("N32").Value = -1573.75

$ws.Range("H74").Value = 2839.8462
$ws.Range("I74").Value = 2872.6287
$ws.Range("J74").Value = 2553
$ws.Range("K74").Value = 2872.6287
$ws.Range("L74").Value = 2553
$ws.Range("M74").Value = -1998.6287
$ws.Range("N74").Value = -4301

$ws.Range("H77").Value = 2839.8462
$ws.Range("I77").Value = 2872.6287
$ws.Range("J77").Value = 2553
$ws.Range("K77").Value = 14363.1435
$ws.Range("L77").Value = 12765
$ws.Range("M77").Value = -9995.143500000002
$ws.Range("N77").Value = -21501

$ws.Range("H122").Value = 14498039
$ws.Range("I122").Value = 33336594
$ws.Range("J122").Value = 6841.5386
$ws.Range("K122").Value = 100009782
$ws.Range("L122").Value = 20524.6158
$ws.Range("M122").Value = -100007332
$ws.Range("N122").Value = -25424.6158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 366.33334
$ws.Range("I14").Value = 249.5
$ws.Range("K14").Value = 249.5
$ws.Range("M14").Value = -79.5

$ws.Range("H31").Value = 2270.8628
$ws.Range("J31").Value = 2919.6155
$ws.Range("L31").Value = 2919.6155
$ws.Range("N31").Value = -3509.6155

$ws.Range("H34").Value = 2270.8628
$ws.Range("J34").Value = 2919.6155
$ws.Range("L34").Value = 2919.6155
$ws.Range("N34").Value = -3323.6155

$ws.Range("H58").Value = 2309.24
$ws.Range("I58").Value = 1548.7
$ws.Range("J58").Value = 2816.2666
$ws.Range("K58").Value = 1548.7
$ws.Range("L58").Value = 2816.2666
$ws.Range("M58").Value = -1345.7
$ws.Range("N58").Value = -3222.2666

$ws.Range("H132").Value = 3041.4583
$ws.Range("I132").Value = 2060.6875
$ws.Range("K132").Value = 6182.0625
$ws.Range("M132").Value = -3652.0625

$ws.Range("H136").Value = 2309.24
$ws.Range("I136").Value = 1548.7
$ws.Range("J136").Value = 2816.2666
$ws.Range("K136").Value = 4646.1
$ws.Range("L136").Value = 8448.799800000001
$ws.Range("M136").Value = -2096.1
$ws.Range("N136").Value = -13548.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 191754.81
$ws.Range("I68").Value = 286746.16
$ws.Range("J68").Value = 1772.1428
$ws.Range("K68").Value = 860238.48
$ws.Range("L68").Value = 5316.428400000001
$ws.Range("M68").Value = -859427.48
$ws.Range("N68").Value = -6938.428400000001

$ws.Range("H71").Value = 191754.81
$ws.Range("I71").Value = 286746.16
$ws.Range("J71").Value = 1772.1428
$ws.Range("K71").Value = 2580715.44
$ws.Range("L71").Value = 15949.2852
$ws.Range("M71").Value = -2576659.44
$ws.Range("N71").Value = -24061.2852

$ws.Range("H107").Value = 1409.8334
$ws.Range("J107").Value = 1300
$ws.Range("L107").Value = 3900
$ws.Range("N107").Value = -7740

$ws.Range("H113").Value = 1194.125
$ws.Range("I113").Value = 361
$ws.Range("J113").Value = 1471.8334
$ws.Range("K113").Value = 1083
$ws.Range("L113").Value = 4415.5002
$ws.Range("M113").Value = 1087
$ws.Range("N113").Value = -8755.5002

$ws.Range("H121").Value = 110893.4
$ws.Range("J121").Value = 184560.67
$ws.Range("L121").Value = 553682.01
$ws.Range("N121").Value = -556302.01

$ws.Range("H131").Value = 8644
$ws.Range("I131").Value = 1351
$ws.Range("K131").Value = 4053
$ws.Range("M131").Value = 987

$ws.Range("H133").Value = 4000
$ws.Range("I133").Value = 4000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 12000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -6940
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3549.4546
$ws.Range("I80").Value = 2888.5
$ws.Range("J80").Value = 3696.3333
$ws.Range("K80").Value = 2888.5
$ws.Range("L80").Value = 3696.3333
$ws.Range("M80").Value = -1890.5
$ws.Range("N80").Value = -5692.3333

$ws.Range("H83").Value = 3549.4546
$ws.Range("I83").Value = 2888.5
$ws.Range("J83").Value = 3696.3333
$ws.Range("K83").Value = 14442.5
$ws.Range("L83").Value = 18481.6665
$ws.Range("M83").Value = -9450.5
$ws.Range("N83").Value = -28465.6665

$ws.Range("H122").Value = 38465116
$ws.Range("I122").Value = 3060
$ws.Range("J122").Value = 55559364
$ws.Range("K122").Value = 9180
$ws.Range("L122").Value = 166678092
$ws.Range("M122").Value = -6730
$ws.Range("N122").Value = -166682992

$ws.Range("H132").Value = 306790.25
$ws.Range("I132").Value = 478817.9
$ws.Range("K132").Value = 1436453.7
$ws.Range("M132").Value = -1433923.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27781000
$ws.Range("I7").Value = 41669228
$ws.Range("K7").Value = 41669228
$ws.Range("M7").Value = -41669116

$ws.Range("H16").Value = 14200
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 17000
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = -2830
$ws.Range("N16").Value = -17340

$ws.Range("H46").Value = 2059.0625
$ws.Range("J46").Value = 2182.2354
$ws.Range("L46").Value = 2182.2354
$ws.Range("N46").Value = -2558.2354

$ws.Range("H93").Value = 3878.8
$ws.Range("I93").Value = 3850
$ws.Range("K93").Value = 3850
$ws.Range("M93").Value = -2602

$ws.Range("H122").Value = 3946.3704
$ws.Range("I122").Value = 2858.1667
$ws.Range("J122").Value = 4816.933
$ws.Range("K122").Value = 8574.500100000001
$ws.Range("L122").Value = 14450.799
$ws.Range("M122").Value = -6124.500100000001
$ws.Range("N122").Value = -19350.799

$ws.Range("H126").Value = 27781000
$ws.Range("I126").Value = 41669228
$ws.Range("K126").Value = 125007684
$ws.Range("M126").Value = -125005214

$ws.Range("H136").Value = 6320.6665
$ws.Range("I136").Value = 3572.6667
$ws.Range("K136").Value = 10718.0001
$ws.Range("M136").Value = -8168.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 530810.1
$ws.Range("I132").Value = 804291.8
$ws.Range("K132").Value = 2412875.4
$ws.Range("M132").Value = -2410345.4

$ws.Range("H136").Value = 7246.2856
$ws.Range("J136").Value = 10621.75
$ws.Range("L136").Value = 31865.25
$ws.Range("N136").Value = -36965.25
